$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy D1's header formatting (bold/border/centered) to the new E1 header cell
# before overwriting D1's text, so both headers keep the matching style.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# Insert new header "EDF-L HOPA" in D1, and push existing "EDF-L GDPA" text to E1
$ws.Range("E1").Value = "EDF-L GDPA"
$ws.Range("D1").Value = "EDF-L HOPA"

# Rows 2-9: column E gets value 50 (column D left unchanged)
$ws.Range("E2").Value = 50
$ws.Range("E3").Value = 50
$ws.Range("E4").Value = 50
$ws.Range("E5").Value = 50
$ws.Range("E6").Value = 50
$ws.Range("E7").Value = 50
$ws.Range("E8").Value = 50
$ws.Range("E9").Value = 50

# Rows 10-21: columns B, C, D updated with new results, column E filled in
$ws.Range("B10").Value = 31
$ws.Range("C10").Value = 47
$ws.Range("D10").Value = 49
$ws.Range("E10").Value = 50

$ws.Range("B11").Value = 28
$ws.Range("C11").Value = 44
$ws.Range("D11").Value = 48
$ws.Range("E11").Value = 50

$ws.Range("B12").Value = 21
$ws.Range("C12").Value = 37
$ws.Range("D12").Value = 45
$ws.Range("E12").Value = 50

$ws.Range("B13").Value = 19
$ws.Range("C13").Value = 36
$ws.Range("D13").Value = 38
$ws.Range("E13").Value = 50

$ws.Range("B14").Value = 12
$ws.Range("C14").Value = 25
$ws.Range("D14").Value = 34
$ws.Range("E14").Value = 49

$ws.Range("B15").Value = 10
$ws.Range("C15").Value = 22
$ws.Range("D15").Value = 30
$ws.Range("E15").Value = 49

$ws.Range("B16").Value = 9
$ws.Range("C16").Value = 17
$ws.Range("D16").Value = 26
$ws.Range("E16").Value = 49

$ws.Range("B17").Value = 6
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 24
$ws.Range("E17").Value = 49

$ws.Range("B18").Value = 5
$ws.Range("C18").Value = 13
$ws.Range("D18").Value = 17
$ws.Range("E18").Value = 43

$ws.Range("B19").Value = 5
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 44

$ws.Range("B20").Value = 4
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = 39

$ws.Range("B21").Value = 2
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 38
